$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in D26:D28
$ws.Range("D26").Value = 338.256
$ws.Range("D27").Value = 389.76
$ws.Range("D28").Value = 478.5
